$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AA2" = 7.5
    "AB2" = 7.5
    "AD2" = 11
    "AI2" = 17
    "AJ2" = 51
    "AK2" = 251
    "G2" = 1.57
    "H2" = 4.2
    "I2" = 5.5
    "J2" = 2.1
    "L2" = 5.5
    "O2" = 1.25
    "P2" = 4
    "Q2" = 1.8
    "R2" = 2
    "U2" = 3
    "V2" = 1.4
    "Y2" = 1.8
    "Z2" = 1.95
    "AB3" = 10
    "AG3" = 8.5
    "AK3" = 301
    "AM3" = 15
    "AP3" = 29
    "AQ3" = 41
    "K3" = 2.05
    "O3" = 1.36
    "P3" = 3.2
    "U3" = 4
    "V3" = 1.25
    "AA4" = 5.5
    "AE4" = 19
    "AL4" = 9.5
    "AR4" = 1.92
    "AS4" = 1.98
    "Q4" = 2.63
    "R4" = 1.5
    "Q5" = 1.9
    "R5" = 1.95
    "AB6" = 12
    "AO6" = 29
    "AP6" = 26
    "G6" = 2.8
    "I6" = 2.75
    "L6" = 3.5
    "AB7" = 5
    "AC7" = 10
    "AD7" = 7
    "AG7" = 8.5
    "AH7" = 10
    "AI7" = 34
    "AJ7" = 151
    "AM7" = 51
    "AN7" = 34
    "AO7" = 151
    "AP7" = 101
    "AQ7" = 101
    "AR7" = 1.58
    "AS7" = 2.39
    "G7" = 1.3
    "H7" = 5
    "I7" = 11
    "J7" = 1.8
    "K7" = 2.3
    "L7" = 11
    "O7" = 1.33
    "P7" = 3.25
    "Q7" = 2.08
    "R7" = 1.73
    "S7" = 3.05
    "T7" = 1.39
    "U7" = 3.75
    "V7" = 1.25
    "Y7" = 2.75
    "Z7" = 1.4
    "W8" = 1.75
    "X8" = 2.05
    "W9" = 1.75
    "X9" = 2.05
    "AA16" = 7.5
    "AB16" = 9
    "AC16" = 8.25
    "AD16" = 15.5
    "AE16" = 14
    "AF16" = 24
    "AG16" = 7.6
    "AH16" = 6.7
    "AK16" = 400
    "AL16" = 12.5
    "AM16" = 24
    "AN16" = 13
    "AO16" = 65
    "AP16" = 35
    "AQ16" = 37
    "G16" = 1.85
    "H16" = 3.45
    "I16" = 4
    "J16" = 2.45
    "K16" = 2.12
    "L16" = 4.3
    "N16" = 7.6
    "O16" = 1.27
    "P16" = 3.45
    "Q16" = 1.8
    "R16" = 1.91
    "U16" = 2.87
    "V16" = 1.37
    "W16" = 1.39
    "X16" = 2.77
    "Y16" = 1.7
    "Z16" = 2.05
    "AA23" = 11
    "AB23" = 19
    "AF23" = 34
    "AG23" = 11
    "AI23" = 13
    "AJ23" = 41
    "AL23" = 8
    "AM23" = 10
    "AO23" = 19
    "AP23" = 17
    "G23" = 3.6
    "I23" = 2
    "J23" = 4
    "L23" = 2.75
    "Q24" = 2
    "R24" = 1.85
    "U24" = 3.5
    "V24" = 1.29
    "AA25" = 5
    "AC25" = 9.5
    "AE25" = 19
    "AG25" = 6.5
    "AI25" = 23
    "AL25" = 9.5
    "AP25" = 51
    "AR25" = 2
    "AS25" = 1.85
    "G25" = 1.76
    "H25" = 3.3
    "K25" = 1.95
    "L25" = 6
    "M25" = 1.1
    "N25" = 7
    "O25" = 1.5
    "P25" = 2.5
    "Q25" = 2.6
    "R25" = 1.48
    "S25" = 4
    "T25" = 1.23
    "U25" = 5.5
    "V25" = 1.14
    "W25" = 1.57
    "X25" = 2.25
    "Y25" = 2.38
    "Z25" = 1.53
    "AA26" = 5
    "AB26" = 8.5
    "AD26" = 21
    "AE26" = 23
    "AH26" = 6.5
    "AI26" = 23
    "AL26" = 7
    "AM26" = 15
    "AN26" = 15
    "AO26" = 41
    "AP26" = 41
    "G26" = 2.15
    "I26" = 3.6
    "J26" = 3.2
    "K26" = 1.83
    "L26" = 4.75
    "Y26" = 2.5
    "Z26" = 1.5
    "G35" = 1.7
    "H35" = 3.5
    "I35" = 5
    "M35" = 1.05
    "N35" = 11
    "Y35" = 1.8
    "Z35" = 1.91
    "AF36" = 29
    "AG36" = 9.5
    "AH36" = 6
    "AI36" = 15
    "AJ36" = 51
    "AK36" = 251
    "AQ36" = 41
    "G36" = 2.1
    "H36" = 3.3
    "K36" = 2.1
    "O36" = 1.33
    "P36" = 3.4
    "Q36" = 2.05
    "R36" = 1.8
    "U36" = 3.5
    "V36" = 1.3
    "W36" = 1.44
    "X36" = 2.63
    "Y36" = 1.8
    "Z36" = 1.91
    "AC48" = 17
    "AG48" = 6.5
    "AL48" = 5
    "AR48" = 2.03
    "AS48" = 1.83
    "G48" = 4.33
    "I48" = 1.83
    "J48" = 5.5
    "K48" = 1.91
    "O48" = 1.53
    "P48" = 2.38
    "Q48" = 2.7
    "R48" = 1.44
    "U48" = 5.5
    "V48" = 1.14
    "Y48" = 2.25
    "Z48" = 1.57
    "AE49" = 19
    "AG49" = 7.5
    "H49" = 3.25
    "O49" = 1.44
    "P49" = 2.63
    "U49" = 4.5
    "V49" = 1.18
    "W49" = 1.53
    "X49" = 2.38
    "Y49" = 2.1
    "Z49" = 1.67
    "O50" = 1.4
    "P50" = 2.75
    "Q50" = 2.25
    "R50" = 1.62
    "U50" = 4.33
    "V50" = 1.2
    "AA51" = 7
    "AB51" = 8
    "AF51" = 26
    "AG51" = 9.5
    "AI51" = 15
    "AK51" = 251
    "AM51" = 23
    "G51" = 1.73
    "I51" = 4.75
    "L51" = 5
    "U51" = 3.25
    "V51" = 1.33
    "W51" = 1.4
    "X51" = 2.75
    "Y51" = 1.91
    "Z51" = 1.91
    "AD54" = 19
    "AG54" = 17
    "AN54" = 12
    "G54" = 2.05
    "H54" = 3.8
    "I54" = 3.25
    "J54" = 2.6
    "N54" = 17
    "Q54" = 1.6
    "R54" = 2.3
    "S54" = 2
    "T54" = 1.85
    "U55" = 4.33
    "V55" = 1.2
    "AD57" = 11
    "AI57" = 23
    "AQ57" = 67
    "AR57" = 1.78
    "AS57" = 2.03
    "G57" = 1.57
    "H57" = 3.3
    "I57" = 6
    "Q57" = 2.35
    "R57" = 1.57
    "Y57" = 2.38
    "Z57" = 1.53
    "O58" = 1.62
    "P58" = 2.2
    "U58" = 6.5
    "V58" = 1.11
    "AR67" = 1.75
    "AS67" = 2.05
    "O67" = 1.4
    "P67" = 2.75
    "Q67" = 2.35
    "R67" = 1.57
    "U67" = 4.33
    "V67" = 1.2
    "AE69" = 13
    "AG69" = 12
    "AI69" = 13
    "AO69" = 51
    "G69" = 1.8
    "H69" = 3.7
    "I69" = 4.5
    "J69" = 2.38
    "M69" = 1.04
    "N69" = 12
    "O69" = 1.22
    "P69" = 4
    "Q69" = 1.8
    "U69" = 2.75
    "V69" = 1.4
    "Y69" = 1.7
    "Z69" = 2.05
    "Q71" = 2.2
    "R71" = 1.65
    "O73" = 1.36
    "P73" = 3
    "Q73" = 2.15
    "R73" = 1.67
    "AC77" = 11
    "AE77" = 23
    "AL77" = 7.5
    "AR77" = 2.1
    "AS77" = 1.78
    "G77" = 2.15
    "H77" = 3
    "I77" = 3.5
    "J77" = 3
    "L77" = 4.33
    "M77" = 1.11
    "O77" = 1.57
    "P77" = 2.25
    "V77" = 1.14
    "Y77" = 2.38
    "Z77" = 1.53
    "I79" = 4.1
    "J79" = 2.88
    "M79" = 1.14
    "N79" = 5.5
    "AA85" = 5.9
    "AB85" = 7.1
    "AC85" = 8
    "AD85" = 12.5
    "AE85" = 14
    "AF85" = 30
    "AH85" = 6.7
    "AI85" = 17
    "AJ85" = 90
    "AK85" = 800
    "AL85" = 13.5
    "AM85" = 32
    "AN85" = 17
    "AO85" = 120
    "AP85" = 60
    "AQ85" = 60
    "G85" = 1.65
    "H85" = 3.4
    "I85" = 5.4
    "J85" = 2.18
    "K85" = 2.1
    "L85" = 5.4
    "O85" = 1.32
    "P85" = 2.82
    "Q85" = 1.93
    "R85" = 1.7
    "U85" = 3.15
    "V85" = 1.26
    "Y85" = 1.88
    "Z85" = 1.72
    "AA86" = 6.2
    "AB86" = 10.5
    "AC86" = 9.25
    "AD86" = 25
    "AJ86" = 80
    "AL86" = 8.25
    "AM86" = 17
    "AN86" = 11.5
    "AO86" = 50
    "G86" = 2.35
    "I86" = 3.3
    "J86" = 3
    "K86" = 1.91
    "L86" = 3.85
    "P86" = 2.42
    "Q86" = 2.25
    "U86" = 3.75
    "Y86" = 1.87
    "Z86" = 1.75
    "AA87" = 7.3
    "AD87" = 14
    "AE87" = 13.5
    "AF87" = 24
    "AG87" = 7.6
    "AH87" = 6.9
    "AI87" = 14.5
    "AJ87" = 65
    "AM87" = 26
    "H87" = 3.5
    "I87" = 4.35
    "J87" = 2.3
    "K87" = 2.15
    "N87" = 7.6
    "W87" = 1.39
    "X87" = 2.8
    "Y87" = 1.75
    "Z87" = 1.98
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
